$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2019-12-31 00:00:00"
$ws.Range("O2").Value = 297962807.98
$ws.Range("P2").Value = 1243.5852653877
$ws.Range("Q2").Value = 2243270706.88
$ws.Range("R2").Value = 9362.572517906599
$ws.Range("S2").Value = 352479236.73
$ws.Range("T2").Value = 1471.1164394113
$ws.Range("U2").Value = -449834818.51
$ws.Range("V2").Value = -1877.4422081394
$ws.Range("W2").ClearContents()
$ws.Range("X2").ClearContents()
$ws.Range("Y2").Value = 455102865.55
$ws.Range("Z2").Value = 1899.4290652265
$ws.Range("AA2").Value = 128694633.65
$ws.Range("AB2").Value = 537.1232444297
$ws.Range("AC2").Value = -23959982.18
$ws.Range("AD2").Value = -126.9793457412
